## TC05_Canine_Filter_Diagnosis-OsteoSarcoma.xlsx
## The "FilesTab" query cell (B4 on the startup sheet) was rewritten:
## the `File Type` column was dropped from the Cypher RETURN clause
## (and the Breed column was dropped too), and the resulting text was
## re-saved, which is what actually happened in Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")
$ws.Activate()

$newFilesQuery = @"

MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
WHERE diag.disease_term IN ['Osteosarcoma']
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS ``File Name``, 
        coalesce(labels(parent)[0], '') AS ``Association``,
        coalesce(f.file_description, '') AS ``Description``,
        coalesce(f.file_format, '') AS ``Format``,
        coalesce(f.file_size, '') AS ``Size``,
        coalesce(c.case_id, '') AS ``Case ID``, 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS ``Study Code``
"@

$ws.Cells.Item(4, 2).Value = $newFilesQuery

# Leave the selection where the user left it after editing the cell.
$ws.Range("B4").Select()
